$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 38885.668
$ws.Range("J3").Value = 38885.668
$ws.Range("L3").Value = 38885.668
$ws.Range("N3").Value = -39113.668

$ws.Range("H53").Value = 26316544
$ws.Range("I53").Value = 71428744
$ws.Range("J53").Value = 1095.0834
$ws.Range("K53").Value = 71428744
$ws.Range("L53").Value = 1095.0834
$ws.Range("M53").Value = -71428107
$ws.Range("N53").Value = -2369.0834

$ws.Range("H69").Value = 4900
$ws.Range("I69").Value = 3166.6667
$ws.Range("J69").Value = 7500
$ws.Range("K69").Value = 9500.000100000001
$ws.Range("L69").Value = 22500
$ws.Range("M69").Value = -8626.000100000001
$ws.Range("N69").Value = -24248

$ws.Range("H72").Value = 4900
$ws.Range("I72").Value = 3166.6667
$ws.Range("J72").Value = 7500
$ws.Range("K72").Value = 28500.0003
$ws.Range("L72").Value = 67500
$ws.Range("M72").Value = -24132.0003
$ws.Range("N72").Value = -76236

$ws.Range("H75").Value = 40000
$ws.Range("J75").Value = 40000
$ws.Range("L75").Value = 40000
$ws.Range("N75").Value = -41872

$ws.Range("H78").Value = 40000
$ws.Range("J78").Value = 40000
$ws.Range("L78").Value = 120000
$ws.Range("N78").Value = -129360

$ws.Range("H102").Value = 38885.668
$ws.Range("J102").Value = 38885.668
$ws.Range("L102").Value = 38885.668
$ws.Range("N102").Value = -45375.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4311708
$ws.Range("I61").Value = 4809074
$ws.Range("J61").Value = 1200
$ws.Range("K61").Value = 4809074
$ws.Range("L61").Value = 1200
$ws.Range("M61").Value = -4808862

$ws.Range("H74").Value = 1105.44
$ws.Range("I74").Value = 1248.9231
$ws.Range("J74").Value = 950
$ws.Range("K74").Value = 1248.9231
$ws.Range("L74").Value = 950
$ws.Range("M74").Value = -374.9231
$ws.Range("N74").Value = -2698

$ws.Range("H77").Value = 1105.44
$ws.Range("I77").Value = 1248.9231
$ws.Range("J77").Value = 950
$ws.Range("K77").Value = 6244.6155
$ws.Range("L77").Value = 4750
$ws.Range("M77").Value = -1876.6155
$ws.Range("N77").Value = -13486

$ws.Range("H102").Value = 2955.6667
$ws.Range("I102").Value = 2711.2856
$ws.Range("K102").Value = 2711.2856
$ws.Range("M102").Value = -1089.2856

$ws.Range("H132").Value = 1401791.9
$ws.Range("I132").Value = 985.40625
$ws.Range("J132").Value = 5884372.5
$ws.Range("K132").Value = 2956.21875
$ws.Range("L132").Value = 17653117.5
$ws.Range("M132").Value = -426.21875
$ws.Range("N132").Value = -17658177.5

$ws.Range("H136").Value = 4311708
$ws.Range("I136").Value = 4809074
$ws.Range("J136").Value = 1200
$ws.Range("K136").Value = 14427222
$ws.Range("L136").Value = 3600
$ws.Range("M136").Value = -14424672

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 19966
$ws.Range("J55").Value = 19966
$ws.Range("L55").Value = 19966
$ws.Range("N55").Value = -20512

$ws.Range("H86").Value = 1165132.1
$ws.Range("I86").Value = 3378
$ws.Range("J86").Value = 2115658.2
$ws.Range("K86").Value = 3378
$ws.Range("L86").Value = 2115658.2
$ws.Range("M86").Value = -2255
$ws.Range("N86").Value = -2117904.2

$ws.Range("H89").Value = 1165132.1
$ws.Range("I89").Value = 3378
$ws.Range("J89").Value = 2115658.2
$ws.Range("K89").Value = 16890
$ws.Range("L89").Value = 10578291
$ws.Range("M89").Value = -11274
$ws.Range("N89").Value = -10589523

$ws.Range("H134").Value = 3973544
$ws.Range("I134").Value = 1744.4783
$ws.Range("J134").Value = 22243822
$ws.Range("K134").Value = 5233.4349
$ws.Range("L134").Value = 66731466
$ws.Range("M134").Value = -2698.4349
$ws.Range("N134").Value = -66736536

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 29900
$ws.Range("I64").Value = 20000
$ws.Range("J64").Value = 39800
$ws.Range("K64").Value = 20000
$ws.Range("L64").Value = 39800
$ws.Range("N64").Value = -40296

$ws.Range("H67").Value = 29900
$ws.Range("I67").Value = 20000
$ws.Range("J67").Value = 39800
$ws.Range("K67").Value = 20000
$ws.Range("L67").Value = 39800
$ws.Range("N67").Value = -41516

$ws.Range("H69").Value = 15000
$ws.Range("J69").Value = 15000
$ws.Range("L69").Value = 15000

$ws.Range("H72").Value = 15000
$ws.Range("J72").Value = 15000
$ws.Range("L72").Value = 45000

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").ClearContents()
$ws.Range("N81").Value = 0

$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").ClearContents()
$ws.Range("N82").Value = 0

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").ClearContents()
$ws.Range("N84").Value = 0

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").ClearContents()
$ws.Range("N85").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 1376.25
$ws.Range("I59").Value = 752.5
$ws.Range("K59").Value = 2257.5
$ws.Range("M59").Value = -1717.5

$ws.Range("H87").Value = 6451.273
$ws.Range("I87").Value = 1988
$ws.Range("J87").Value = 8125
$ws.Range("K87").Value = 5964
$ws.Range("L87").Value = 24375
$ws.Range("M87").Value = -4716
$ws.Range("N87").Value = -26871

$ws.Range("H90").Value = 6451.273
$ws.Range("I90").Value = 1988
$ws.Range("J90").Value = 8125
$ws.Range("K90").Value = 17892
$ws.Range("L90").Value = 73125
$ws.Range("M90").Value = -11652
$ws.Range("N90").Value = -85605

$ws.Range("H96").Value = 3000
$ws.Range("J96").Value = 3000
$ws.Range("L96").Value = 9000

$ws.Range("H98").Value = 467.88235
$ws.Range("J98").Value = 501.1
$ws.Range("L98").Value = 1503.3
$ws.Range("N98").Value = -4499.3

$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H102").Value = 4941.6665
$ws.Range("I102").Value = 4825
$ws.Range("K102").Value = 14475

$ws.Range("H103").Value = 102045
$ws.Range("I103").Value = 612.5
$ws.Range("J103").Value = 169666.67
$ws.Range("K103").Value = 1837.5
$ws.Range("L103").Value = 509000.01
$ws.Range("M103").Value = -958.5
$ws.Range("N103").Value = -510758.01

$ws.Range("H104").Value = 10000
$ws.Range("J104").Value = 10000
$ws.Range("L104").Value = 30000

$ws.Range("H105").Value = 3285.7144
$ws.Range("J105").Value = 3285.7144
$ws.Range("L105").Value = 9857.143199999999
$ws.Range("N105").Value = -15099.1432

$ws.Range("H106").Value = 2000
$ws.Range("J106").Value = 2000
$ws.Range("L106").Value = 6000
$ws.Range("N106").Value = -7892

$ws.Range("H129").Value = 1237.6666
$ws.Range("J129").Value = 1733.2222
$ws.Range("L129").Value = 5199.6666
$ws.Range("N129").Value = -15199.6666

$ws.Range("H131").Value = 878.6900000000001
$ws.Range("I131").Value = 557.5
$ws.Range("J131").Value = 892.07294
$ws.Range("K131").Value = 1672.5
$ws.Range("L131").Value = 2676.21882
$ws.Range("M131").Value = 3367.5
$ws.Range("N131").Value = -12756.21882

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 14168.9
$ws.Range("I132").Value = 4499.4
$ws.Range("J132").Value = 23838.4
$ws.Range("K132").Value = 13498.2
$ws.Range("L132").Value = 71515.20000000001
$ws.Range("M132").Value = -10968.2
$ws.Range("N132").Value = -76575.20000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1431.2667
$ws.Range("I68").Value = 1442.3636
$ws.Range("J68").Value = 1400.75
$ws.Range("K68").Value = 1442.3636
$ws.Range("L68").Value = 1400.75
$ws.Range("M68").Value = -693.3635999999999
$ws.Range("N68").Value = -2898.75

$ws.Range("H71").Value = 1431.2667
$ws.Range("I71").Value = 1442.3636
$ws.Range("J71").Value = 1400.75
$ws.Range("K71").Value = 7211.817999999999
$ws.Range("L71").Value = 7003.75
$ws.Range("M71").Value = -3467.817999999999
$ws.Range("N71").Value = -14491.75

$ws.Range("H93").Value = 1194.2858
$ws.Range("I93").Value = 1052
$ws.Range("J93").Value = 1550
$ws.Range("K93").Value = 1052
$ws.Range("L93").Value = 1550
$ws.Range("M93").Value = 196
$ws.Range("N93").Value = -4046
